$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new "2022-Q4" sheet right after "总计" (i.e. before the
#    sheet that is currently "2022-Q3"), matching the new tab order:
#    总计, 2022-Q4, 2022-Q3, 2022-Q2, 2021-Q1
# ---------------------------------------------------------------------------
$oldQ3Ws = $wb.Worksheets.Item(2)
$q4Ws = $wb.Worksheets.Add($oldQ3Ws)
$q4Ws.Name = "2022-Q4"

# NOTE: worksheet handles in this host track by *position*, so after the
# insert above the variable that used to point at "2022-Q3" now resolves to
# the freshly-added sheet instead. Re-resolve every sheet we still need by
# name/index AFTER the structural change.
$totalWs = $wb.Worksheets.Item("总计")
$q3Ws = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------------
# 2) Populate the new "2022-Q4" sheet with the fund table, re-using the
#    header / index-column formatting (style) already used on the other
#    per-quarter sheets so the look matches exactly.
# ---------------------------------------------------------------------------
$q3Ws.Range("B1:H1").Copy()
$q4Ws.Range("B1:H1").PasteSpecial(-4122)

$q3Ws.Range("A2").Copy()
$q4Ws.Range("A2").PasteSpecial(-4122)
$q4Ws.Range("A3").PasteSpecial(-4122)

$q4Ws.Range("B1").Value = "基金代码"
$q4Ws.Range("C1").Value = "基金名称"
$q4Ws.Range("D1").Value = "基金规模"
$q4Ws.Range("E1").Value = "股票总仓位"
$q4Ws.Range("F1").Value = "仓位占比"
$q4Ws.Range("G1").Value = "持有市值(亿元)"
$q4Ws.Range("H1").Value = "仓位排名"

$q4Ws.Range("A2").Value = 0
$q4Ws.Range("C2").Value = "金鹰民安回报一年定期开放混合A"
$q4Ws.Range("H2").Value = 1

$q4Ws.Range("A3").Value = 1
$q4Ws.Range("C3").Value = "金鹰民安回报一年定期开放混合C"
$q4Ws.Range("H3").Value = 1

# Columns that hold digit-strings (fund code, size, position, weight, value)
# must stay TEXT, not become numbers (which would also eat leading zeros in
# the fund codes) - write them as a formula returning the literal string and
# then convert that formula to a static value, which keeps the text type
# without tagging the cell with a "stored as text" quote-prefix style.
$q4Ws.Range("B2").Value = '="006972"'
$q4Ws.Range("D2").Value = '="21.93"'
$q4Ws.Range("E2").Value = '="39.73"'
$q4Ws.Range("F2").Value = '="2.39"'
$q4Ws.Range("G2").Value = '="0.5241"'

$q4Ws.Range("B3").Value = '="007735"'
$q4Ws.Range("D3").Value = '="2.16"'
$q4Ws.Range("E3").Value = '="39.73"'
$q4Ws.Range("F3").Value = '="2.39"'
$q4Ws.Range("G3").Value = '="0.0516"'

$q4Ws.Range("B2:G3").Copy()
$q4Ws.Range("B2:G3").PasteSpecial(-4163)

# ---------------------------------------------------------------------------
# 3) Update the "总计" summary sheet: the new 2022-Q4 row is inserted at the
#    top of the data (row 2) and every following row shifts down by one, so
#    rewrite all four data rows with their new contents.
# ---------------------------------------------------------------------------
$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q4"
$totalWs.Range("C2").Value = 2
$totalWs.Range("D2").Value = 0.58

$totalWs.Range("A3").Value = 1
$totalWs.Range("B3").Value = "2022-Q3"
$totalWs.Range("C3").Value = 6
$totalWs.Range("D3").Value = 0.35

$totalWs.Range("A4").Value = 2
$totalWs.Range("B4").Value = "2022-Q2"
$totalWs.Range("C4").Value = 2
$totalWs.Range("D4").Value = 0.06

$totalWs.Range("A4").Copy()
$totalWs.Range("A5").PasteSpecial(-4122)
$totalWs.Range("A5").Value = 3
$totalWs.Range("B5").Value = "2021-Q1"
$totalWs.Range("C5").Value = 2
$totalWs.Range("D5").Value = 0.08
